# Apply cryptos list update (values from the diff) cell by cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.422.68'
$ws.Range('E2').Value = '  -2.57%  '
$ws.Range('D3').Value = '3.692.19'
$ws.Range('E3').Value = '  -3.07%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '692.35'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.27'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -5.41%  '
$ws.Range('D7').Value = '3.690.69'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -4.82%  '
$ws.Range('E10').Value = '  -8.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.38'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.441'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -5.52%  '
$ws.Range('E13').Value = '  -5.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.29'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -7.44%  '
$ws.Range('D15').Value = '4.313.65'
$ws.Range('E15').Value = '  -3.06%  '
$ws.Range('D16').Value = '3.699.36'
$ws.Range('E16').Value = '  -3.54%  '
$ws.Range('D17').Value = '69.483.13'
$ws.Range('E17').Value = '  -2.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.13'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -7.73%  '
$ws.Range('E20').Value = '  -8.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '478.78'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -6.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.01'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.72%  '
$ws.Range('E23').Value = '  -7.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.97'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.82%  '
$ws.Range('D25').Value = '3.836.99'
$ws.Range('E26').Value = '  -9.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.34'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.48'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -8.52%  '
$ws.Range('E30').Value = '  -11.10%  '
$ws.Range('E31').Value = '  -10.10%  '
$ws.Range('E32').Value = '  -7.46%  '
$ws.Range('E33').Value = '  -7.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.168'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.91'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -7.30%  '
$ws.Range('D37').Value = '3.658.44'
$ws.Range('E37').Value = '  -2.96%  '
$ws.Range('E38').Value = '  -7.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.27'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.35'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('E41').Value = '  -8.27%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  -6.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '163.48'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -6.05%  '
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '30.02'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.37%  '
$ws.Range('E48').Value = '  -14.99%  '
$ws.Range('E49').Value = '  -0.46%  '
$ws.Range('E50').Value = '  -2.41%  '
$ws.Range('E51').Value = '  -9.22%  '
